$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Garanti nepřednášejí" fix: guarantor (department-level) rows had been
# mixed up with the individually-taught rows. Swap course-name / zkratka
# pairs so each row again reflects the correct course <-> department code.

# Rows 8 & 9 ("Počítačové modelování I" / "Programování A", K107 / K103)
$ws.Range("A8").Value = "Programování A"
$ws.Range("B8").Value = "K103"
$ws.Range("A9").Value = "Počítačové modelování I"
$ws.Range("B9").Value = "K107"

# Rows 27 & 29 (KSPP / KRSPP) - row 28 (RSPP) is untouched
$ws.Range("B27").Value = "KRSPP"
$ws.Range("B29").Value = "KSPP"
